$d = $word.ActiveDocument

# The paragraphs that were justified ("both") get their alignment cleared
# back to the (unset/default) value, i.e.
#   <w:pPr><w:jc w:val="both"/></w:pPr>  ->  <w:pPr/>
# Paragraphs aligned center/right (or already default) are left untouched.
#   wdAlignParagraphLeft = 0, wdAlignParagraphJustify = 3
$wdAlignParagraphLeft = 0
$wdAlignParagraphJustify = 3

For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    If ($para.Format.Alignment -eq $wdAlignParagraphJustify) {
        $para.Format.Alignment = $wdAlignParagraphLeft
    }
}
